# Web 120 / Homework 10 & 11
# Adds two new homework score columns (H10 -> M, H11 -> N) with student
# scores, flips a few "extra point" flags in the summary table, and moves
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new column headers for Homework 10 and Homework 11 ---
$ws.Range("M1").Value = "H10"
$ws.Range("N1").Value = "H11"

# --- Homework 10 / 11 scores for each student (rows 2-16) ---
# Rows 2-10 already use the centred 2-decimal style, so a plain value
# assignment keeps the existing formatting.
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 10

$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 9

$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 10

$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 10

$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 8.5

$ws.Range("M7").Value = 9
$ws.Range("N7").Value = 9

$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0

$ws.Range("M9").Value = 10
$ws.Range("N9").Value = 10.5

$ws.Range("M10").Value = 10
$ws.Range("N10").Value = 7

# Rows 11-16 were still on the left-aligned "empty" style, so line them up
# with the rest of the column (centre aligned, like column L) before
# filling in the scores.
$ws.Range("M11:N16").HorizontalAlignment = -4108

$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 8.5

$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 11

$ws.Range("M13").Value = 10
$ws.Range("N13").Value = 8.5

$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 5

$ws.Range("M15").Value = 10
$ws.Range("N15").Value = 8

$ws.Range("M16").Value = 10
$ws.Range("N16").Value = 8.5

# --- Summary table: grant the "extra point" flag (column C) to the
#     students who now qualify because of the new homework scores. ---
$ws.Range("C22").Value = 1
$ws.Range("C24").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C32").Value = 1
$ws.Range("C33").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("C36").Value = 1

# --- Move the active selection to where the editor left off ---
$ws.Range("M26").Select()
